$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 155.57143
$ws.Range("I2").Value = 145.8
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 145.8
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = -32.80000000000001
$ws.Range("N2").Value = -406
# Row 11
$ws.Range("H11").Value = 1018.1
$ws.Range("I11").Value = 1018.1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1018.1
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -878.1
# Row 31
$ws.Range("H31").Value = 3329.6
$ws.Range("I31").Value = 3329.6
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 9988.799999999999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -9758.799999999999
# Row 38
$ws.Range("H38").Value = 1938.6875
$ws.Range("I38").Value = 1001.46155
$ws.Range("J38").Value = 6000
$ws.Range("K38").Value = 3004.38465
$ws.Range("L38").Value = 18000
$ws.Range("M38").Value = -2632.38465
$ws.Range("N38").Value = -18744
# Row 39
$ws.Range("H39").Value = 5240.2856
$ws.Range("I39").Value = 1295
$ws.Range("J39").Value = 9185.571
$ws.Range("K39").Value = 3885
$ws.Range("L39").Value = 27556.713
$ws.Range("M39").Value = -3589
# Row 42
$ws.Range("H42").Value = 4122.5
$ws.Range("I42").Value = 5598
$ws.Range("J42").Value = 3384.75
$ws.Range("K42").Value = 16794
$ws.Range("L42").Value = 10154.25
$ws.Range("M42").Value = -16564
$ws.Range("N42").Value = -10614.25
# Row 53
$ws.Range("H53").Value = 1151.25
$ws.Range("I53").Value = 222.25
$ws.Range("J53").Value = 3009.25
$ws.Range("K53").Value = 222.25
$ws.Range("L53").Value = 3009.25
$ws.Range("M53").Value = 414.75
$ws.Range("N53").Value = -4283.25
# Row 132
$ws.Range("H132").Value = 2554.75
$ws.Range("I132").Value = 2630.8572
$ws.Range("J132").Value = 2326.4285
$ws.Range("K132").Value = 7892.571599999999
$ws.Range("L132").Value = 6979.2855
$ws.Range("M132").Value = -5362.571599999999
$ws.Range("N132").Value = -12039.2855
# Row 137
$ws.Range("H137").Value = 37039492
$ws.Range("I137").Value = 52633340
$ws.Range("J137").Value = 4105.25
$ws.Range("K137").Value = 157900020
$ws.Range("L137").Value = 12315.75
$ws.Range("M137").Value = -157897470
$ws.Range("N137").Value = -17415.75
# Row 138
$ws.Range("H138").Value = 3538.3193
$ws.Range("I138").Value = 1791.6072
$ws.Range("J138").Value = 4649.864
$ws.Range("K138").Value = 5374.821599999999
$ws.Range("L138").Value = 13949.592
$ws.Range("M138").Value = -234.8215999999993
$ws.Range("N138").Value = -24229.592

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 546250.3
$ws.Range("I2").Value = 775335.8
$ws.Range("J2").Value = 2172.125
$ws.Range("K2").Value = 775335.8
$ws.Range("L2").Value = 2172.125
$ws.Range("M2").Value = -775222.8
$ws.Range("N2").Value = -2398.125
# Row 45
$ws.Range("H45").Value = 1053.8334
$ws.Range("I45").Value = 1031
$ws.Range("J45").Value = 1099.5
$ws.Range("K45").Value = 1031
$ws.Range("L45").Value = 1099.5
$ws.Range("M45").Value = -654
# Row 74
$ws.Range("H74").Value = 34487292
$ws.Range("I74").Value = 38465692
$ws.Range("J74").Value = 7833.6665
$ws.Range("K74").Value = 38465692
$ws.Range("L74").Value = 7833.6665
$ws.Range("M74").Value = -38464818
# Row 77
$ws.Range("H77").Value = 34487292
$ws.Range("I77").Value = 38465692
$ws.Range("J77").Value = 7833.6665
$ws.Range("K77").Value = 192328460
$ws.Range("L77").Value = 39168.3325
$ws.Range("M77").Value = -192324092
# Row 97
$ws.Range("H97").Value = 175.5
$ws.Range("I97").Value = 150.66667
$ws.Range("J97").Value = 250
$ws.Range("K97").Value = 150.66667
$ws.Range("L97").Value = 250
$ws.Range("M97").Value = 345.33333
# Row 102
$ws.Range("H102").Value = 6668042
$ws.Range("I102").Value = 7693695
$ws.Range("J102").Value = 1299.5
$ws.Range("K102").Value = 7693695
$ws.Range("L102").Value = 1299.5
$ws.Range("M102").Value = -7692073
# Row 116
$ws.Range("H116").Value = 546250.3
$ws.Range("I116").Value = 775335.8
$ws.Range("J116").Value = 2172.125
$ws.Range("K116").Value = 775335.8
$ws.Range("L116").Value = 2172.125
$ws.Range("M116").Value = -773041.8
$ws.Range("N116").Value = -6760.125
# Row 132
$ws.Range("H132").Value = 3033910.2
$ws.Range("I132").Value = 3128576
$ws.Range("J132").Value = 4605
$ws.Range("K132").Value = 9385728
$ws.Range("L132").Value = 13815
$ws.Range("M132").Value = -9383198

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 546250.3
$ws.Range("I3").Value = 775335.8
$ws.Range("J3").Value = 2172.125
$ws.Range("K3").Value = 775335.8
$ws.Range("L3").Value = 2172.125
$ws.Range("M3").Value = -775221.8
$ws.Range("N3").Value = -2400.125
# Row 11
$ws.Range("H11").Value = 255.5
$ws.Range("I11").Value = 226.6
$ws.Range("J11").Value = 276.14285
$ws.Range("K11").Value = 226.6
$ws.Range("L11").Value = 276.14285
$ws.Range("M11").Value = -86.59999999999999
$ws.Range("N11").Value = -556.14285
# Row 20
$ws.Range("H20").Value = 1189.3334
$ws.Range("I20").Value = 1071.8572
$ws.Range("J20").Value = 1292.125
$ws.Range("K20").Value = 1071.8572
$ws.Range("L20").Value = 1292.125
$ws.Range("M20").Value = -824.8571999999999
$ws.Range("N20").Value = -1786.125
# Row 107
$ws.Range("H107").Value = 72839.86
$ws.Range("I107").Value = 1555.3636
$ws.Range("J107").Value = 334216.34
$ws.Range("K107").Value = 1555.3636
$ws.Range("L107").Value = 334216.34
$ws.Range("M107").Value = 364.6364000000001
# Row 134
$ws.Range("H134").Value = 12954420
$ws.Range("I134").Value = 13295200
$ws.Range("J134").Value = 4791
$ws.Range("K134").Value = 39885600
$ws.Range("L134").Value = 14373
$ws.Range("M134").Value = -39883065

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8799.538
$ws.Range("I31").Value = 5212.125
$ws.Range("J31").Value = 14539.4
$ws.Range("K31").Value = 5212.125
$ws.Range("L31").Value = 14539.4
$ws.Range("M31").Value = -4917.125
$ws.Range("N31").Value = -15129.4
# Row 32
$ws.Range("H32").Value = 24460
$ws.Range("I32").Value = 13380
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 13380
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -13064
# Row 34
$ws.Range("H34").Value = 8799.538
$ws.Range("I34").Value = 5212.125
$ws.Range("J34").Value = 14539.4
$ws.Range("K34").Value = 5212.125
$ws.Range("L34").Value = 14539.4
$ws.Range("M34").Value = -5010.125
$ws.Range("N34").Value = -14943.4
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 134
$ws.Range("H134").Value = 38891736
$ws.Range("I134").Value = 22730250
$ws.Range("J134").Value = 83335820
$ws.Range("K134").Value = 68190750
$ws.Range("L134").Value = 250007460
$ws.Range("M134").Value = -68188215
$ws.Range("N134").Value = -250012530

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 72.94118
$ws.Range("I2").Value = 57.1
$ws.Range("J2").Value = 95.57143000000001
$ws.Range("K2").Value = 342.6
$ws.Range("L2").Value = 573.42858
$ws.Range("M2").Value = -229.6
$ws.Range("N2").Value = -799.42858
# Row 94
$ws.Range("H94").Value = 16638.363
$ws.Range("I94").Value = 4799
$ws.Range("J94").Value = 26504.5
$ws.Range("K94").Value = 14397
$ws.Range("L94").Value = 79513.5
$ws.Range("M94").Value = -13721
$ws.Range("N94").Value = -80865.5
# Row 121
$ws.Range("H121").Value = 170672
$ws.Range("I121").Value = 339999.66
$ws.Range("J121").Value = 1344.3334
$ws.Range("K121").Value = 1019998.98
$ws.Range("L121").Value = 4033.0002
$ws.Range("M121").Value = -1018688.98
$ws.Range("N121").Value = -6653.0002
# Row 122
$ws.Range("H122").Value = 699.3333
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 299
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 2691
$ws.Range("M122").Value = -11050
# Row 129
$ws.Range("H129").Value = 3823.5
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 4165.85
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 12497.55
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -22497.55
# Row 131
$ws.Range("H131").Value = 1179.8889
$ws.Range("I131").Value = 1014.875
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = 3044.625
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = 1995.375
$ws.Range("N131").Value = -17580
# Row 137
$ws.Range("H137").Value = 9093284
# Row 139
$ws.Range("H139").Value = 683.8889
$ws.Range("I139").Value = 581.875
$ws.Range("J139").Value = 1500
$ws.Range("K139").Value = 1745.625
$ws.Range("L139").Value = 4500
$ws.Range("M139").Value = 3394.375
$ws.Range("N139").Value = -14780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7486.625
$ws.Range("I70").Value = 7598.8
$ws.Range("J70").Value = 7299.6665
$ws.Range("K70").Value = 7598.8
$ws.Range("L70").Value = 7299.6665
$ws.Range("M70").Value = -7328.8
# Row 73
$ws.Range("H73").Value = 7486.625
$ws.Range("I73").Value = 7598.8
$ws.Range("J73").Value = 7299.6665
$ws.Range("K73").Value = 7598.8
$ws.Range("L73").Value = 7299.6665
$ws.Range("M73").Value = -6662.8
# Row 97
$ws.Range("H97").Value = 1354
$ws.Range("I97").Value = 1334.7
$ws.Range("J97").Value = 1431.2
$ws.Range("K97").Value = 1334.7
$ws.Range("L97").Value = 1431.2
$ws.Range("M97").Value = -838.7
$ws.Range("N97").Value = -2423.2
# Row 126
$ws.Range("H126").Value = 3932.2666
$ws.Range("I126").Value = 3666.8333
$ws.Range("J126").Value = 4994
$ws.Range("K126").Value = 11000.4999
$ws.Range("L126").Value = 14982
$ws.Range("M126").Value = -8530.499899999999
# Row 132
$ws.Range("H132").Value = 9037604
$ws.Range("I132").Value = 9858875
$ws.Range("J132").Value = 3623.3333
$ws.Range("K132").Value = 29576625
$ws.Range("L132").Value = 10869.9999
$ws.Range("M132").Value = -29574095

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 62502684
$ws.Range("I22").Value = 2855.375
$ws.Range("J22").Value = 125002510
$ws.Range("K22").Value = 2855.375
$ws.Range("L22").Value = 125002510
$ws.Range("M22").Value = -2560.375
$ws.Range("N22").Value = -125003100
# Row 27
$ws.Range("H27").Value = 62502684
$ws.Range("I27").Value = 2855.375
$ws.Range("J27").Value = 125002510
$ws.Range("K27").Value = 2855.375
$ws.Range("L27").Value = 125002510
$ws.Range("M27").Value = -2748.375
$ws.Range("N27").Value = -125002724
# Row 46
$ws.Range("H46").Value = 2825
$ws.Range("I46").Value = 2728.5715
$ws.Range("J46").Value = 3500
$ws.Range("K46").Value = 2728.5715
$ws.Range("L46").Value = 3500
$ws.Range("M46").Value = -2540.5715
# Row 55
$ws.Range("H55").Value = 545.6
$ws.Range("I55").Value = 308.375
$ws.Range("J55").Value = 816.7143
$ws.Range("K55").Value = 308.375
$ws.Range("L55").Value = 816.7143
$ws.Range("M55").Value = -135.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 19492.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 19492.5
$ws.Range("N45").Value = -20474.5
$ws.Range("M45").ClearContents()
# Row 107
$ws.Range("H107").Value = 1147.125
$ws.Range("I107").Value = 879.5
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 2638.5
$ws.Range("L107").Value = 5850
$ws.Range("M107").Value = -718.5
$ws.Range("N107").Value = -9690
# Row 122
$ws.Range("H122").Value = 3040.2666
$ws.Range("I122").Value = 2969.5386
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8908.6158
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6458.6158
